$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new risk entry on row 14 (Risk ID 15)
$ws.Range("B14").Value = "7/18/2017"
$ws.Range("C14").Value = "Problemas Técnicos"
$ws.Range("D14").Value = "Atualmente só um integrante possui notebook, em caso de problema, ficamos incapazes de apresentar."
$ws.Range("E14").Value = "Technical"
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 1
$ws.Range("I14").Value = "Todos"
$ws.Range("J14").Value = "Tirar prints das telas que desenvolvemos, para em caso de problemas, ser possível mostrar o andamento do projeto."

# Adjust row height to fit the new content
$ws.Rows.Item(14).RowHeight = 43.5

# Update sheet view: zoom level and active selection
$excel.ActiveWindow.Zoom = 70
$ws.Range("C15").Select() | Out-Null
